# Update "paises.xlsx" country COVID figures with the latest snapshot (data pulled
# at 18 Junio 2020 17:55) and re-sort the country table descending by "Casos totales",
# matching the upstream worldometer refresh used for this tracker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$n = 216
$data = New-Object 'object[,]' $n,8

$data[0,0] = "Estados Unidos"
$data[0,1] = 2243406
$data[0,2] = 8935
$data[0,3] = 919140
$data[0,4] = 1204171
$data[0,5] = 0
$data[0,6] = 154
$data[0,7] = 120095

$data[1,0] = "Brasil"
$data[1,1] = 960309
$data[1,2] = 0
$data[1,3] = 503507
$data[1,4] = 410137
$data[1,5] = 0
$data[1,6] = 0
$data[1,7] = 46665

$data[2,0] = "Rusia"
$data[2,1] = 561091
$data[2,2] = 7790
$data[2,3] = 313963
$data[2,4] = 239468
$data[2,5] = 0
$data[2,6] = 182
$data[2,7] = 7660

$data[3,0] = "India"
$data[3,1] = 371734
$data[3,2] = 4470
$data[3,3] = 196894
$data[3,4] = 162480
$data[3,5] = 0
$data[3,6] = 98
$data[3,7] = 12360

$data[4,0] = "Reino Unido"
$data[4,1] = 300469
$data[4,2] = 1218
$data[4,3] = 0
$data[4,4] = 0
$data[4,5] = 0
$data[4,6] = 135
$data[4,7] = 42288

$data[5,0] = "España"
$data[5,1] = 291763
$data[5,2] = 0
$data[5,3] = 0
$data[5,4] = 0
$data[5,5] = 0
$data[5,6] = 0
$data[5,7] = 27136

$data[6,0] = "Peru"
$data[6,1] = 240908
$data[6,2] = 0
$data[6,3] = 128622
$data[6,4] = 105029
$data[6,5] = 0
$data[6,6] = 0
$data[6,7] = 7257

$data[7,0] = "Italia"
$data[7,1] = 237828
$data[7,2] = 0
$data[7,3] = 179455
$data[7,4] = 23925
$data[7,5] = 0
$data[7,6] = 0
$data[7,7] = 34448

$data[8,0] = "Chile"
$data[8,1] = 220628
$data[8,2] = 0
$data[8,3] = 156232
$data[8,4] = 60781
$data[8,5] = 0
$data[8,6] = 0
$data[8,7] = 3615

$data[9,0] = "Iran"
$data[9,1] = 197647
$data[9,2] = 2596
$data[9,3] = 156991
$data[9,4] = 31384
$data[9,5] = 0
$data[9,6] = 87
$data[9,7] = 9272

$data[10,0] = "Alemania"
$data[10,1] = 189504
$data[10,2] = 0
$data[10,3] = 173600
$data[10,4] = 6977
$data[10,5] = 0
$data[10,6] = 0
$data[10,7] = 8927

$data[11,0] = "Turquia"
$data[11,1] = 182727
$data[11,2] = 0
$data[11,3] = 154640
$data[11,4] = 23226
$data[11,5] = 0
$data[11,6] = 0
$data[11,7] = 4861

$data[12,0] = "Pakistan"
$data[12,1] = 160118
$data[12,2] = 5358
$data[12,3] = 59215
$data[12,4] = 97810
$data[12,5] = 0
$data[12,6] = 118
$data[12,7] = 3093

$data[13,0] = "Mexico"
$data[13,1] = 159793
$data[13,2] = 4930
$data[13,3] = 119355
$data[13,4] = 21358
$data[13,5] = 0
$data[13,6] = 770
$data[13,7] = 19080

$data[14,0] = "Francia"
$data[14,1] = 158174
$data[14,2] = 0
$data[14,3] = 73667
$data[14,4] = 54932
$data[14,5] = 0
$data[14,6] = 0
$data[14,7] = 29575

$data[15,0] = "Arabia Saudita"
$data[15,1] = 145991
$data[15,2] = 4757
$data[15,3] = 93915
$data[15,4] = 50937
$data[15,5] = 0
$data[15,6] = 48
$data[15,7] = 1139

$data[16,0] = "Banglades"
$data[16,1] = 102292
$data[16,2] = 3803
$data[16,3] = 40164
$data[16,4] = 60785
$data[16,5] = 0
$data[16,6] = 38
$data[16,7] = 1343

$data[17,0] = "Canada"
$data[17,1] = 99853
$data[17,2] = 0
$data[17,3] = 62017
$data[17,4] = 29582
$data[17,5] = 0
$data[17,6] = 0
$data[17,7] = 8254

$data[18,0] = "Catar"
$data[18,1] = 84441
$data[18,2] = 1267
$data[18,3] = 63642
$data[18,4] = 20713
$data[18,5] = 0
$data[18,6] = 4
$data[18,7] = 86

$data[19,0] = "China"
$data[19,1] = 83293
$data[19,2] = 28
$data[19,3] = 78394
$data[19,4] = 265
$data[19,5] = 0
$data[19,6] = 0
$data[19,7] = 4634

$data[20,0] = "Sudafrica"
$data[20,1] = 80412
$data[20,2] = 0
$data[20,3] = 44331
$data[20,4] = 34407
$data[20,5] = 0
$data[20,6] = 0
$data[20,7] = 1674

$data[21,0] = "Belgica"
$data[21,1] = 60348
$data[21,2] = 104
$data[21,3] = 16724
$data[21,4] = 33941
$data[21,5] = 0
$data[21,6] = 8
$data[21,7] = 9683

$data[22,0] = "Colombia"
$data[22,1] = 57046
$data[22,2] = 0
$data[22,3] = 21326
$data[22,4] = 33856
$data[22,5] = 0
$data[22,6] = 0
$data[22,7] = 1864

$data[23,0] = "Bielorrusia"
$data[23,1] = 56657
$data[23,2] = 625
$data[23,3] = 34023
$data[23,4] = 22303
$data[23,5] = 0
$data[23,6] = 7
$data[23,7] = 331

$data[24,0] = "Suecia"
$data[24,1] = 56043
$data[24,2] = 1481
$data[24,3] = 0
$data[24,4] = 0
$data[24,5] = 0
$data[24,6] = 12
$data[24,7] = 5053

$data[25,0] = "Paises Bajos"
$data[25,1] = 49319
$data[25,2] = 115
$data[25,3] = 0
$data[25,4] = 0
$data[25,5] = 0
$data[25,6] = 4
$data[25,7] = 6078

$data[26,0] = "Egipto"
$data[26,1] = 49219
$data[26,2] = 0
$data[26,3] = 13141
$data[26,4] = 34228
$data[26,5] = 0
$data[26,6] = 0
$data[26,7] = 1850

$data[27,0] = "Ecuador"
$data[27,1] = 48490
$data[27,2] = 0
$data[27,3] = 23881
$data[27,4] = 20602
$data[27,5] = 0
$data[27,6] = 0
$data[27,7] = 4007

$data[28,0] = "Emiratos Arabes Unidos"
$data[28,1] = 43752
$data[28,2] = 388
$data[28,3] = 30241
$data[28,4] = 13213
$data[28,5] = 0
$data[28,6] = 3
$data[28,7] = 298

$data[29,0] = "Indonesia"
$data[29,1] = 42762
$data[29,2] = 1331
$data[29,3] = 16798
$data[29,4] = 23625
$data[29,5] = 0
$data[29,6] = 63
$data[29,7] = 2339

$data[30,0] = "Singapur"
$data[30,1] = 41473
$data[30,2] = 257
$data[30,3] = 32712
$data[30,4] = 8735
$data[30,5] = 0
$data[30,6] = 0
$data[30,7] = 26

$data[31,0] = "Portugal"
$data[31,1] = 38089
$data[31,2] = 417
$data[31,3] = 24010
$data[31,4] = 12555
$data[31,5] = 0
$data[31,6] = 1
$data[31,7] = 1524

$data[32,0] = "Kuwait"
$data[32,1] = 38074
$data[32,2] = 541
$data[32,3] = 29512
$data[32,4] = 8254
$data[32,5] = 0
$data[32,6] = 2
$data[32,7] = 308

$data[33,0] = "Argentina"
$data[33,1] = 35552
$data[33,2] = 0
$data[33,3] = 10721
$data[33,4] = 23902
$data[33,5] = 0
$data[33,6] = 16
$data[33,7] = 929

$data[34,0] = "Ucrania"
$data[34,1] = 34063
$data[34,2] = 829
$data[34,3] = 15447
$data[34,4] = 17650
$data[34,5] = 0
$data[34,6] = 23
$data[34,7] = 966

$data[35,0] = "Suiza"
$data[35,1] = 31200
$data[35,2] = 13
$data[35,3] = 28900
$data[35,4] = 344
$data[35,5] = 0
$data[35,6] = 0
$data[35,7] = 1956

$data[36,0] = "Polonia"
$data[36,1] = 31015
$data[36,2] = 314
$data[36,3] = 15317
$data[36,4] = 14382
$data[36,5] = 0
$data[36,6] = 30
$data[36,7] = 1316

$data[37,0] = "Filipinas"
$data[37,1] = 27799
$data[37,2] = 561
$data[37,3] = 7090
$data[37,4] = 19593
$data[37,5] = 0
$data[37,6] = 8
$data[37,7] = 1116

$data[38,0] = "Afganistan"
$data[38,1] = 27532
$data[38,2] = 658
$data[38,3] = 7660
$data[38,4] = 19326
$data[38,5] = 0
$data[38,6] = 42
$data[38,7] = 546

$data[39,0] = "Oman"
$data[39,1] = 26818
$data[39,2] = 739
$data[39,3] = 13264
$data[39,4] = 13435
$data[39,5] = 0
$data[39,6] = 3
$data[39,7] = 119

$data[40,0] = "Irlanda"
$data[40,1] = 25341
$data[40,2] = 0
$data[40,3] = 22698
$data[40,4] = 933
$data[40,5] = 0
$data[40,6] = 0
$data[40,7] = 1710

$data[41,0] = "Republica Dominicana"
$data[41,1] = 24645
$data[41,2] = 540
$data[41,3] = 14293
$data[41,4] = 9717
$data[41,5] = 0
$data[41,6] = 2
$data[41,7] = 635

$data[42,0] = "Irak"
$data[42,1] = 24254
$data[42,2] = 0
$data[42,3] = 10770
$data[42,4] = 12711
$data[42,5] = 0
$data[42,6] = 0
$data[42,7] = 773

$data[43,0] = "Rumania"
$data[43,1] = 23080
$data[43,2] = 320
$data[43,3] = 16308
$data[43,4] = 5299
$data[43,5] = 0
$data[43,6] = 22
$data[43,7] = 1473

$data[44,0] = "Panama"
$data[44,1] = 22597
$data[44,2] = 0
$data[44,3] = 13774
$data[44,4] = 8353
$data[44,5] = 0
$data[44,6] = 0
$data[44,7] = 470

$data[45,0] = "Bolivia"
$data[45,1] = 20685
$data[45,2] = 802
$data[45,3] = 4002
$data[45,4] = 16004
$data[45,5] = 0
$data[45,6] = 20
$data[45,7] = 679

$data[46,0] = "Barein"
$data[46,1] = 19961
$data[46,2] = 0
$data[46,3] = 14185
$data[46,4] = 5723
$data[46,5] = 0
$data[46,6] = 4
$data[46,7] = 53

$data[47,0] = "Israel"
$data[47,1] = 19894
$data[47,2] = 111
$data[47,3] = 15499
$data[47,4] = 4092
$data[47,5] = 0
$data[47,6] = 0
$data[47,7] = 303

$data[48,0] = "Armenia"
$data[48,1] = 18698
$data[48,2] = 665
$data[48,3] = 7560
$data[48,4] = 10829
$data[48,5] = 0
$data[48,6] = 7
$data[48,7] = 309

$data[49,0] = "Nigeria"
$data[49,1] = 17735
$data[49,2] = 0
$data[49,3] = 5967
$data[49,4] = 11299
$data[49,5] = 0
$data[49,6] = 0
$data[49,7] = 469

$data[50,0] = "Japon"
$data[50,1] = 17628
$data[50,2] = 0
$data[50,3] = 15850
$data[50,4] = 847
$data[50,5] = 0
$data[50,6] = 0
$data[50,7] = 931

$data[51,0] = "Austria"
$data[51,1] = 17223
$data[51,2] = 20
$data[51,3] = 16101
$data[51,4] = 434
$data[51,5] = 0
$data[51,6] = 1
$data[51,7] = 688

$data[52,0] = "Kazajistan"
$data[52,1] = 15877
$data[52,2] = 335
$data[52,3] = 10065
$data[52,4] = 5712
$data[52,5] = 0
$data[52,6] = 3
$data[52,7] = 100

$data[53,0] = "Moldavia"
$data[53,1] = 13106
$data[53,2] = 374
$data[53,3] = 7252
$data[53,4] = 5410
$data[53,5] = 0
$data[53,6] = 11
$data[53,7] = 444

$data[54,0] = "Ghana"
$data[54,1] = 12929
$data[54,2] = 339
$data[54,3] = 4468
$data[54,4] = 8395
$data[54,5] = 0
$data[54,6] = 0
$data[54,7] = 66

$data[55,0] = "Serbia"
$data[55,1] = 12616
$data[55,2] = 94
$data[55,3] = 11511
$data[55,4] = 847
$data[55,5] = 0
$data[55,6] = 1
$data[55,7] = 258

$data[56,0] = "Dinamarca"
$data[56,1] = 12344
$data[56,2] = 50
$data[56,3] = 11242
$data[56,4] = 502
$data[56,5] = 0
$data[56,6] = 2
$data[56,7] = 600

$data[57,0] = "Corea del Sur"
$data[57,1] = 12257
$data[57,2] = 59
$data[57,3] = 10800
$data[57,4] = 1177
$data[57,5] = 0
$data[57,6] = 1
$data[57,7] = 280

$data[58,0] = "Azerbaiyan"
$data[58,1] = 11329
$data[58,2] = 338
$data[58,3] = 6192
$data[58,4] = 4998
$data[58,5] = 0
$data[58,6] = 6
$data[58,7] = 139

$data[59,0] = "Argelia"
$data[59,1] = 11268
$data[59,2] = 0
$data[59,3] = 7943
$data[59,4] = 2526
$data[59,5] = 0
$data[59,6] = 0
$data[59,7] = 799

$data[60,0] = "Guatemala"
$data[60,1] = 11251
$data[60,2] = 545
$data[60,3] = 2200
$data[60,4] = 8619
$data[60,5] = 0
$data[60,6] = 14
$data[60,7] = 432

$data[61,0] = "Honduras"
$data[61,1] = 10299
$data[61,2] = 643
$data[61,3] = 1116
$data[61,4] = 8847
$data[61,5] = 0
$data[61,6] = 6
$data[61,7] = 336

$data[62,0] = "Chequia"
$data[62,1] = 10176
$data[62,2] = 14
$data[62,3] = 7401
$data[62,4] = 2442
$data[62,5] = 0
$data[62,6] = 0
$data[62,7] = 333

$data[63,0] = "Camerun"
$data[63,1] = 9864
$data[63,2] = 0
$data[63,3] = 5570
$data[63,4] = 4018
$data[63,5] = 0
$data[63,6] = 0
$data[63,7] = 276

$data[64,0] = "Marruecos"
$data[64,1] = 9042
$data[64,2] = 45
$data[64,3] = 7999
$data[64,4] = 830
$data[64,5] = 0
$data[64,6] = 0
$data[64,7] = 213

$data[65,0] = "Noruega"
$data[65,1] = 8707
$data[65,2] = 15
$data[65,3] = 8138
$data[65,4] = 325
$data[65,5] = 0
$data[65,6] = 1
$data[65,7] = 244

$data[66,0] = "Malasia"
$data[66,1] = 8529
$data[66,2] = 14
$data[66,3] = 8000
$data[66,4] = 408
$data[66,5] = 0
$data[66,6] = 0
$data[66,7] = 121

$data[67,0] = "Sudan"
$data[67,1] = 8020
$data[67,2] = 0
$data[67,3] = 2966
$data[67,4] = 4567
$data[67,5] = 0
$data[67,6] = 0
$data[67,7] = 487

$data[68,0] = "Nepal"
$data[68,1] = 7848
$data[68,2] = 671
$data[68,3] = 1186
$data[68,4] = 6640
$data[68,5] = 0
$data[68,6] = 2
$data[68,7] = 22

$data[69,0] = "Australia"
$data[69,1] = 7391
$data[69,2] = 21
$data[69,3] = 6877
$data[69,4] = 412
$data[69,5] = 0
$data[69,6] = 0
$data[69,7] = 102

$data[70,0] = "Finlandia"
$data[70,1] = 7119
$data[70,2] = 2
$data[70,3] = 6200
$data[70,4] = 593
$data[70,5] = 0
$data[70,6] = 0
$data[70,7] = 326

$data[71,0] = "Costa de Marfil"
$data[71,1] = 6063
$data[71,2] = 0
$data[71,3] = 2749
$data[71,4] = 3266
$data[71,5] = 0
$data[71,6] = 0
$data[71,7] = 48

$data[72,0] = "Uzbekistan"
$data[72,1] = 5730
$data[72,2] = 48
$data[72,3] = 4166
$data[72,4] = 1545
$data[72,5] = 0
$data[72,6] = 0
$data[72,7] = 19

$data[73,0] = "Senegal"
$data[73,1] = 5475
$data[73,2] = 106
$data[73,3] = 3716
$data[73,4] = 1683
$data[73,5] = 0
$data[73,6] = 3
$data[73,7] = 76

$data[74,0] = "Consejo Danes para los Refugiados"
$data[74,1] = 5283
$data[74,2] = 183
$data[74,3] = 685
$data[74,4] = 4481
$data[74,5] = 0
$data[74,6] = 2
$data[74,7] = 117

$data[75,0] = "Tayikistan"
$data[75,1] = 5279
$data[75,2] = 58
$data[75,3] = 3762
$data[75,4] = 1466
$data[75,5] = 0
$data[75,6] = 0
$data[75,7] = 51

$data[76,0] = "Haiti"
$data[76,1] = 4688
$data[76,2] = 141
$data[76,3] = 24
$data[76,4] = 4582
$data[76,5] = 0
$data[76,6] = 2
$data[76,7] = 82

$data[77,0] = "Guinea"
$data[77,1] = 4668
$data[77,2] = 0
$data[77,3] = 3364
$data[77,4] = 1278
$data[77,5] = 0
$data[77,6] = 0
$data[77,7] = 26

$data[78,0] = "Republica de Macedonia"
$data[78,1] = 4664
$data[78,2] = 182
$data[78,3] = 1836
$data[78,4] = 2612
$data[78,5] = 0
$data[78,6] = 6
$data[78,7] = 216

$data[79,0] = "Republica de Yibuti"
$data[79,1] = 4545
$data[79,2] = 0
$data[79,3] = 3411
$data[79,4] = 1091
$data[79,5] = 0
$data[79,6] = 0
$data[79,7] = 43

$data[80,0] = "Kenia"
$data[80,1] = 4257
$data[80,2] = 213
$data[80,3] = 1459
$data[80,4] = 2681
$data[80,5] = 0
$data[80,6] = 10
$data[80,7] = 117

$data[81,0] = "Gabon"
$data[81,1] = 4229
$data[81,2] = 0
$data[81,3] = 1505
$data[81,4] = 2694
$data[81,5] = 0
$data[81,6] = 0
$data[81,7] = 30

$data[82,0] = "El Salvador"
$data[82,1] = 4200
$data[82,2] = 134
$data[82,3] = 2235
$data[82,4] = 1883
$data[82,5] = 0
$data[82,6] = 3
$data[82,7] = 82

$data[83,0] = "Luxemburgo"
$data[83,1] = 4091
$data[83,2] = 6
$data[83,3] = 3940
$data[83,4] = 41
$data[83,5] = 0
$data[83,6] = 0
$data[83,7] = 110

$data[84,0] = "Hungria"
$data[84,1] = 4079
$data[84,2] = 1
$data[84,3] = 2564
$data[84,4] = 947
$data[84,5] = 0
$data[84,6] = 1
$data[84,7] = 568

$data[85,0] = "Etiopia"
$data[85,1] = 3954
$data[85,2] = 195
$data[85,3] = 934
$data[85,4] = 2955
$data[85,5] = 0
$data[85,6] = 2
$data[85,7] = 65

$data[86,0] = "Bulgaria"
$data[86,1] = 3542
$data[86,2] = 0
$data[86,3] = 1880
$data[86,4] = 1478
$data[86,5] = 0
$data[86,6] = 0
$data[86,7] = 184

$data[87,0] = "Venezuela"
$data[87,1] = 3386
$data[87,2] = 0
$data[87,3] = 835
$data[87,4] = 2523
$data[87,5] = 0
$data[87,6] = 0
$data[87,7] = 28

$data[88,0] = "Grecia"
$data[88,1] = 3227
$data[88,2] = 24
$data[88,3] = 1374
$data[88,4] = 1665
$data[88,5] = 0
$data[88,6] = 1
$data[88,7] = 188

$data[89,0] = "Bosnia y Herzegovina"
$data[89,1] = 3174
$data[89,2] = 33
$data[89,3] = 2219
$data[89,4] = 787
$data[89,5] = 0
$data[89,6] = 0
$data[89,7] = 168

$data[90,0] = "Tailandia"
$data[90,1] = 3141
$data[90,2] = 6
$data[90,3] = 2997
$data[90,4] = 86
$data[90,5] = 0
$data[90,6] = 0
$data[90,7] = 58

$data[91,0] = "Somalia"
$data[91,1] = 2719
$data[91,2] = 23
$data[91,3] = 724
$data[91,4] = 1907
$data[91,5] = 0
$data[91,6] = 0
$data[91,7] = 88

$data[92,0] = "Kirguistan"
$data[92,1] = 2657
$data[92,2] = 95
$data[92,3] = 1933
$data[92,4] = 693
$data[92,5] = 0
$data[92,6] = 1
$data[92,7] = 31

$data[93,0] = "Republica de Africa Central"
$data[93,1] = 2564
$data[93,2] = 0
$data[93,3] = 402
$data[93,4] = 2144
$data[93,5] = 0
$data[93,6] = 0
$data[93,7] = 18

$data[94,0] = "Mayotte"
$data[94,1] = 2383
$data[94,2] = 38
$data[94,3] = 2066
$data[94,4] = 288
$data[94,5] = 0
$data[94,6] = 0
$data[94,7] = 29

$data[95,0] = "Cuba"
$data[95,1] = 2295
$data[95,2] = 15
$data[95,3] = 2020
$data[95,4] = 190
$data[95,5] = 0
$data[95,6] = 1
$data[95,7] = 85

$data[96,0] = "Croacia"
$data[96,1] = 2269
$data[96,2] = 11
$data[96,3] = 2142
$data[96,4] = 20
$data[96,5] = 0
$data[96,6] = 0
$data[96,7] = 107

$data[97,0] = "Mauritania"
$data[97,1] = 2223
$data[97,2] = 0
$data[97,3] = 427
$data[97,4] = 1701
$data[97,5] = 0
$data[97,6] = 0
$data[97,7] = 95

$data[98,0] = "Maldivas"
$data[98,1] = 2120
$data[98,2] = 0
$data[98,3] = 1677
$data[98,4] = 435
$data[98,5] = 0
$data[98,6] = 0
$data[98,7] = 8

$data[99,0] = "Estonia"
$data[99,1] = 1977
$data[99,2] = 0
$data[99,3] = 1748
$data[99,4] = 160
$data[99,5] = 0
$data[99,6] = 0
$data[99,7] = 69

$data[100,0] = "Sri Lanka"
$data[100,1] = 1926
$data[100,2] = 2
$data[100,3] = 1421
$data[100,4] = 494
$data[100,5] = 0
$data[100,6] = 0
$data[100,7] = 11

$data[101,0] = "Mali"
$data[101,1] = 1906
$data[101,2] = 16
$data[101,3] = 1192
$data[101,4] = 607
$data[101,5] = 0
$data[101,6] = 0
$data[101,7] = 107

$data[102,0] = "Costa Rica"
$data[102,1] = 1871
$data[102,2] = 0
$data[102,3] = 899
$data[102,4] = 960
$data[102,5] = 0
$data[102,6] = 0
$data[102,7] = 12

$data[103,0] = "Nicaragua"
$data[103,1] = 1823
$data[103,2] = 0
$data[103,3] = 1238
$data[103,4] = 521
$data[103,5] = 0
$data[103,6] = 0
$data[103,7] = 64

$data[104,0] = "Islandia"
$data[104,1] = 1816
$data[104,2] = 1
$data[104,3] = 1801
$data[104,4] = 5
$data[104,5] = 0
$data[104,6] = 0
$data[104,7] = 10

$data[105,0] = "Sudan del Sur"
$data[105,1] = 1813
$data[105,2] = 0
$data[105,3] = 89
$data[105,4] = 1693
$data[105,5] = 0
$data[105,6] = 0
$data[105,7] = 31

$data[106,0] = "Albania"
$data[106,1] = 1788
$data[106,2] = 66
$data[106,3] = 1086
$data[106,4] = 663
$data[106,5] = 0
$data[106,6] = 1
$data[106,7] = 39

$data[107,0] = "Lituania"
$data[107,1] = 1784
$data[107,2] = 6
$data[107,3] = 1449
$data[107,4] = 259
$data[107,5] = 0
$data[107,6] = 0
$data[107,7] = 76

$data[108,0] = "Guayana Francesa"
$data[108,1] = 1758
$data[108,2] = 204
$data[108,3] = 784
$data[108,4] = 969
$data[108,5] = 0
$data[108,6] = 0
$data[108,7] = 5

$data[109,0] = "Guinea Ecuatorial"
$data[109,1] = 1664
$data[109,2] = 0
$data[109,3] = 515
$data[109,4] = 1117
$data[109,5] = 0
$data[109,6] = 0
$data[109,7] = 32

$data[110,0] = "Eslovaquia"
$data[110,1] = 1562
$data[110,2] = 1
$data[110,3] = 1443
$data[110,4] = 91
$data[110,5] = 0
$data[110,6] = 0
$data[110,7] = 28

$data[111,0] = "Eslovenia"
$data[111,1] = 1511
$data[111,2] = 8
$data[111,3] = 1359
$data[111,4] = 43
$data[111,5] = 0
$data[111,6] = 0
$data[111,7] = 109

$data[112,0] = "Nueva Zelanda"
$data[112,1] = 1507
$data[112,2] = 1
$data[112,3] = 1482
$data[112,4] = 3
$data[112,5] = 0
$data[112,6] = 0
$data[112,7] = 22

$data[113,0] = "Libano"
$data[113,1] = 1495
$data[113,2] = 6
$data[113,3] = 944
$data[113,4] = 519
$data[113,5] = 0
$data[113,6] = 0
$data[113,7] = 32

$data[114,0] = "Guinea-Bisau"
$data[114,1] = 1492
$data[114,2] = 0
$data[114,3] = 153
$data[114,4] = 1324
$data[114,5] = 0
$data[114,6] = 0
$data[114,7] = 15

$data[115,0] = "Zambia"
$data[115,1] = 1416
$data[115,2] = 4
$data[115,3] = 1144
$data[115,4] = 261
$data[115,5] = 0
$data[115,6] = 0
$data[115,7] = 11

$data[116,0] = "Madagascar"
$data[116,1] = 1403
$data[116,2] = 25
$data[116,3] = 463
$data[116,4] = 927
$data[116,5] = 0
$data[116,6] = 1
$data[116,7] = 13

$data[117,0] = "Paraguay"
$data[117,1] = 1308
$data[117,2] = 0
$data[117,3] = 711
$data[117,4] = 584
$data[117,5] = 0
$data[117,6] = 0
$data[117,7] = 13

$data[118,0] = "Sierra Leona"
$data[118,1] = 1249
$data[118,2] = 0
$data[118,3] = 707
$data[118,4] = 491
$data[118,5] = 0
$data[118,6] = 0
$data[118,7] = 51

$data[119,0] = "Tunez"
$data[119,1] = 1132
$data[119,2] = 4
$data[119,3] = 1006
$data[119,4] = 76
$data[119,5] = 0
$data[119,6] = 0
$data[119,7] = 50

$data[120,0] = "Hong Kong"
$data[120,1] = 1125
$data[120,2] = 4
$data[120,3] = 1072
$data[120,4] = 49
$data[120,5] = 0
$data[120,6] = 0
$data[120,7] = 4

$data[121,0] = "Letonia"
$data[121,1] = 1108
$data[121,2] = 4
$data[121,3] = 903
$data[121,4] = 175
$data[121,5] = 0
$data[121,6] = 0
$data[121,7] = 30

$data[122,0] = "Niger"
$data[122,1] = 1020
$data[122,2] = 0
$data[122,3] = 893
$data[122,4] = 60
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 67

$data[123,0] = "Jordania"
$data[123,1] = 1001
$data[123,2] = 14
$data[123,3] = 697
$data[123,4] = 295
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 9

$data[124,0] = "Republica de Chipre"
$data[124,1] = 985
$data[124,2] = 0
$data[124,3] = 816
$data[124,4] = 151
$data[124,5] = 0
$data[124,6] = 0
$data[124,7] = 18

$data[125,0] = "Yemen"
$data[125,1] = 902
$data[125,2] = 0
$data[125,3] = 271
$data[125,4] = 387
$data[125,5] = 0
$data[125,6] = 0
$data[125,7] = 244

$data[126,0] = "Burkina Faso"
$data[126,1] = 899
$data[126,2] = 0
$data[126,3] = 810
$data[126,4] = 36
$data[126,5] = 0
$data[126,6] = 0
$data[126,7] = 53

$data[127,0] = "Georgia"
$data[127,1] = 893
$data[127,2] = 5
$data[127,3] = 739
$data[127,4] = 140
$data[127,5] = 0
$data[127,6] = 0
$data[127,7] = 14

$data[128,0] = "Congo"
$data[128,1] = 883
$data[128,2] = 0
$data[128,3] = 391
$data[128,4] = 465
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 27

$data[129,0] = "Principado de Andorra"
$data[129,1] = 855
$data[129,2] = 1
$data[129,3] = 792
$data[129,4] = 11
$data[129,5] = 0
$data[129,6] = 0
$data[129,7] = 52

$data[130,0] = "Republica del Chad"
$data[130,1] = 854
$data[130,2] = 0
$data[130,3] = 733
$data[130,4] = 47
$data[130,5] = 0
$data[130,6] = 0
$data[130,7] = 74

$data[131,0] = "Uruguay"
$data[131,1] = 849
$data[131,2] = 0
$data[131,3] = 810
$data[131,4] = 15
$data[131,5] = 0
$data[131,6] = 0
$data[131,7] = 24

$data[132,0] = "Cabo Verde"
$data[132,1] = 823
$data[132,2] = 31
$data[132,3] = 377
$data[132,4] = 439
$data[132,5] = 0
$data[132,6] = 0
$data[132,7] = 7

$data[133,0] = "Uganda"
$data[133,1] = 741
$data[133,2] = 9
$data[133,3] = 486
$data[133,4] = 255
$data[133,5] = 0
$data[133,6] = 0
$data[133,7] = 0

$data[134,0] = "Crucero"
$data[134,1] = 712
$data[134,2] = 0
$data[134,3] = 651
$data[134,4] = 48
$data[134,5] = 0
$data[134,6] = 0
$data[134,7] = 13

$data[135,0] = "San Marino"
$data[135,1] = 696
$data[135,2] = 0
$data[135,3] = 609
$data[135,4] = 42
$data[135,5] = 0
$data[135,6] = 3
$data[135,7] = 45

$data[136,0] = "Santo Tome y Principe"
$data[136,1] = 683
$data[136,2] = 0
$data[136,3] = 188
$data[136,4] = 483
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 12

$data[137,0] = "Malta"
$data[137,1] = 663
$data[137,2] = 1
$data[137,3] = 610
$data[137,4] = 44
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 9

$data[138,0] = "Mozambique"
$data[138,1] = 662
$data[138,2] = 11
$data[138,3] = 175
$data[138,4] = 483
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 4

$data[139,0] = "Ruanda"
$data[139,1] = 639
$data[139,2] = 0
$data[139,3] = 347
$data[139,4] = 290
$data[139,5] = 0
$data[139,6] = 0
$data[139,7] = 2

$data[140,0] = "Jamaica"
$data[140,1] = 626
$data[140,2] = 5
$data[140,3] = 451
$data[140,4] = 165
$data[140,5] = 0
$data[140,6] = 0
$data[140,7] = 10

$data[141,0] = "Benin"
$data[141,1] = 597
$data[141,2] = 25
$data[141,3] = 238
$data[141,4] = 348
$data[141,5] = 0
$data[141,6] = 2
$data[141,7] = 11

$data[142,0] = "Estado de Palestina"
$data[142,1] = 579
$data[142,2] = 24
$data[142,3] = 415
$data[142,4] = 161
$data[142,5] = 0
$data[142,6] = 0
$data[142,7] = 3

$data[143,0] = "Malaui"
$data[143,1] = 572
$data[143,2] = 0
$data[143,3] = 73
$data[143,4] = 493
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 6

$data[144,0] = "Suazilandia"
$data[144,1] = 563
$data[144,2] = 0
$data[144,3] = 262
$data[144,4] = 297
$data[144,5] = 0
$data[144,6] = 0
$data[144,7] = 4

$data[145,0] = "Togo"
$data[145,1] = 544
$data[145,2] = 0
$data[145,3] = 353
$data[145,4] = 178
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 13

$data[146,0] = "Liberia"
$data[146,1] = 542
$data[146,2] = 26
$data[146,3] = 250
$data[146,4] = 259
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 33

$data[147,0] = "Tanzania"
$data[147,1] = 509
$data[147,2] = 0
$data[147,3] = 183
$data[147,4] = 305
$data[147,5] = 0
$data[147,6] = 0
$data[147,7] = 21

$data[148,0] = "Reunion"
$data[148,1] = 502
$data[148,2] = 5
$data[148,3] = 460
$data[148,4] = 41
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 1

$data[149,0] = "Libia"
$data[149,1] = 500
$data[149,2] = 0
$data[149,3] = 78
$data[149,4] = 412
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 10

$data[150,0] = "Taiwan"
$data[150,1] = 446
$data[150,2] = 1
$data[150,3] = 434
$data[150,4] = 5
$data[150,5] = 0
$data[150,6] = 0
$data[150,7] = 7

$data[151,0] = "Zimbabue"
$data[151,1] = 401
$data[151,2] = 0
$data[151,3] = 63
$data[151,4] = 334
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 4

$data[152,0] = "Vietnam"
$data[152,1] = 342
$data[152,2] = 7
$data[152,3] = 325
$data[152,4] = 17
$data[152,5] = 0
$data[152,6] = 0
$data[152,7] = 0

$data[153,0] = "Mauricio"
$data[153,1] = 337
$data[153,2] = 0
$data[153,3] = 325
$data[153,4] = 2
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 10

$data[154,0] = "Isla de Man"
$data[154,1] = 336
$data[154,2] = 0
$data[154,3] = 312
$data[154,4] = 0
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 24

$data[155,0] = "Montenegro"
$data[155,1] = 333
$data[155,2] = 0
$data[155,3] = 315
$data[155,4] = 9
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 9

$data[156,0] = "Birmania"
$data[156,1] = 263
$data[156,2] = 1
$data[156,3] = 185
$data[156,4] = 72
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 6

$data[157,0] = "Surinam"
$data[157,1] = 261
$data[157,2] = 0
$data[157,3] = 48
$data[157,4] = 207
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 6

$data[158,0] = "Martinica"
$data[158,1] = 221
$data[158,2] = 19
$data[158,3] = 98
$data[158,4] = 109
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 14

$data[159,0] = "Mongolia"
$data[159,1] = 201
$data[159,2] = 4
$data[159,3] = 127
$data[159,4] = 74
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 0

$data[160,0] = "Comoras"
$data[160,1] = 197
$data[160,2] = 0
$data[160,3] = 127
$data[160,4] = 67
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 3

$data[161,0] = "Islas Caimanes"
$data[161,1] = 193
$data[161,2] = 0
$data[161,3] = 132
$data[161,4] = 60
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 1

$data[162,0] = "Siria"
$data[162,1] = 187
$data[162,2] = 9
$data[162,3] = 78
$data[162,4] = 102
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 7

$data[163,0] = "Islas Feroe"
$data[163,1] = 187
$data[163,2] = 0
$data[163,3] = 187
$data[163,4] = 0
$data[163,5] = 0
$data[163,6] = 0
$data[163,7] = 0

$data[164,0] = "Gibraltar"
$data[164,1] = 176
$data[164,2] = 0
$data[164,3] = 176
$data[164,4] = 0
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 0

$data[165,0] = "Guyana"
$data[165,1] = 171
$data[165,2] = 0
$data[165,3] = 102
$data[165,4] = 57
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 12

$data[166,0] = "Guadalupe"
$data[166,1] = 171
$data[166,2] = 0
$data[166,3] = 157
$data[166,4] = 0
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 14

$data[167,0] = "Angola"
$data[167,1] = 155
$data[167,2] = 0
$data[167,3] = 64
$data[167,4] = 84
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 7

$data[168,0] = "Bermudas"
$data[168,1] = 144
$data[168,2] = 0
$data[168,3] = 128
$data[168,4] = 7
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 9

$data[169,0] = "Brunei"
$data[169,1] = 141
$data[169,2] = 0
$data[169,3] = 138
$data[169,4] = 0
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 3

$data[170,0] = "Eritrea"
$data[170,1] = 131
$data[170,2] = 0
$data[170,3] = 39
$data[170,4] = 92
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 0

$data[171,0] = "Camboya"
$data[171,1] = 129
$data[171,2] = 1
$data[171,3] = 126
$data[171,4] = 3
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 0

$data[172,0] = "Trinidad yTobago"
$data[172,1] = 123
$data[172,2] = 0
$data[172,3] = 109
$data[172,4] = 6
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 8

$data[173,0] = "Burundi"
$data[173,1] = 104
$data[173,2] = 0
$data[173,3] = 75
$data[173,4] = 28
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 1

$data[174,0] = "Bahamas"
$data[174,1] = 104
$data[174,2] = 0
$data[174,3] = 72
$data[174,4] = 21
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 11

$data[175,0] = "Aruba"
$data[175,1] = 101
$data[175,2] = 0
$data[175,3] = 98
$data[175,4] = 0
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 3

$data[176,0] = "Monaco"
$data[176,1] = 99
$data[176,2] = 0
$data[176,3] = 94
$data[176,4] = 1
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 4

$data[177,0] = "Barbados"
$data[177,1] = 97
$data[177,2] = 0
$data[177,3] = 85
$data[177,4] = 5
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 7

$data[178,0] = "Liechtenstein"
$data[178,1] = 82
$data[178,2] = 0
$data[178,3] = 55
$data[178,4] = 26
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 1

$data[179,0] = "Botsuana"
$data[179,1] = 79
$data[179,2] = 0
$data[179,3] = 25
$data[179,4] = 53
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 1

$data[180,0] = "San Martin (Parte Holandesa)"
$data[180,1] = 77
$data[180,2] = 0
$data[180,3] = 62
$data[180,4] = 0
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 15

$data[181,0] = "Butan"
$data[181,1] = 67
$data[181,2] = 0
$data[181,3] = 25
$data[181,4] = 42
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 0

$data[182,0] = "Polinesia Francesa"
$data[182,1] = 60
$data[182,2] = 0
$data[182,3] = 60
$data[182,4] = 0
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 0

$data[183,0] = "Macao"
$data[183,1] = 45
$data[183,2] = 0
$data[183,3] = 45
$data[183,4] = 0
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 0

$data[184,0] = "San Martin (Parte Francesa)"
$data[184,1] = 42
$data[184,2] = 0
$data[184,3] = 36
$data[184,4] = 3
$data[184,5] = 0
$data[184,6] = 0
$data[184,7] = 3

$data[185,0] = "Puerto Rico"
$data[185,1] = 39
$data[185,2] = 0
$data[185,3] = 1
$data[185,4] = 36
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 2

$data[186,0] = "Namibia"
$data[186,1] = 39
$data[186,2] = 3
$data[186,3] = 19
$data[186,4] = 20
$data[186,5] = 0
$data[186,6] = 0
$data[186,7] = 0

$data[187,0] = "Gambia"
$data[187,1] = 34
$data[187,2] = 0
$data[187,3] = 24
$data[187,4] = 9
$data[187,5] = 0
$data[187,6] = 0
$data[187,7] = 1

$data[188,0] = "Guam"
$data[188,1] = 32
$data[188,2] = 0
$data[188,3] = 0
$data[188,4] = 31
$data[188,5] = 0
$data[188,6] = 0
$data[188,7] = 1

$data[189,0] = "San Vicente y las Granadinas"
$data[189,1] = 29
$data[189,2] = 0
$data[189,3] = 25
$data[189,4] = 4
$data[189,5] = 0
$data[189,6] = 0
$data[189,7] = 0

$data[190,0] = "Antigua y Barbuda"
$data[190,1] = 26
$data[190,2] = 0
$data[190,3] = 22
$data[190,4] = 1
$data[190,5] = 0
$data[190,6] = 0
$data[190,7] = 3

$data[191,0] = "Timor Oriental"
$data[191,1] = 24
$data[191,2] = 0
$data[191,3] = 24
$data[191,4] = 0
$data[191,5] = 0
$data[191,6] = 0
$data[191,7] = 0

$data[192,0] = "Curazao"
$data[192,1] = 23
$data[192,2] = 0
$data[192,3] = 19
$data[192,4] = 3
$data[192,5] = 0
$data[192,6] = 0
$data[192,7] = 1

$data[193,0] = "Granada"
$data[193,1] = 23
$data[193,2] = 0
$data[193,3] = 22
$data[193,4] = 1
$data[193,5] = 0
$data[193,6] = 0
$data[193,7] = 0

$data[194,0] = "Belice"
$data[194,1] = 22
$data[194,2] = 0
$data[194,3] = 16
$data[194,4] = 4
$data[194,5] = 0
$data[194,6] = 0
$data[194,7] = 2

$data[195,0] = "Nueva Caledonia"
$data[195,1] = 21
$data[195,2] = 0
$data[195,3] = 21
$data[195,4] = 0
$data[195,5] = 0
$data[195,6] = 0
$data[195,7] = 0

$data[196,0] = "Santa Lucia"
$data[196,1] = 19
$data[196,2] = 0
$data[196,3] = 18
$data[196,4] = 1
$data[196,5] = 0
$data[196,6] = 0
$data[196,7] = 0

$data[197,0] = "Laos"
$data[197,1] = 19
$data[197,2] = 0
$data[197,3] = 19
$data[197,4] = 0
$data[197,5] = 0
$data[197,6] = 0
$data[197,7] = 0

$data[198,0] = "Dominica"
$data[198,1] = 18
$data[198,2] = 0
$data[198,3] = 18
$data[198,4] = 0
$data[198,5] = 0
$data[198,6] = 0
$data[198,7] = 0

$data[199,0] = "Fiyi"
$data[199,1] = 18
$data[199,2] = 0
$data[199,3] = 18
$data[199,4] = 0
$data[199,5] = 0
$data[199,6] = 0
$data[199,7] = 0

$data[200,0] = "Islas Virgenes de los Estados Unidos"
$data[200,1] = 17
$data[200,2] = 0
$data[200,3] = 0
$data[200,4] = 17
$data[200,5] = 0
$data[200,6] = 0
$data[200,7] = 0

$data[201,0] = "San Cristobal y Nieves"
$data[201,1] = 15
$data[201,2] = 0
$data[201,3] = 15
$data[201,4] = 0
$data[201,5] = 0
$data[201,6] = 0
$data[201,7] = 0

$data[202,0] = "Islas Malvinas"
$data[202,1] = 13
$data[202,2] = 0
$data[202,3] = 13
$data[202,4] = 0
$data[202,5] = 0
$data[202,6] = 0
$data[202,7] = 0

$data[203,0] = "Groenlandia"
$data[203,1] = 13
$data[203,2] = 0
$data[203,3] = 13
$data[203,4] = 0
$data[203,5] = 0
$data[203,6] = 0
$data[203,7] = 0

$data[204,0] = "Santa Sede"
$data[204,1] = 12
$data[204,2] = 0
$data[204,3] = 12
$data[204,4] = 0
$data[204,5] = 0
$data[204,6] = 0
$data[204,7] = 0

$data[205,0] = "Islas Turcas y Caicos"
$data[205,1] = 12
$data[205,2] = 0
$data[205,3] = 11
$data[205,4] = 0
$data[205,5] = 0
$data[205,6] = 0
$data[205,7] = 1

$data[206,0] = "Montserrat"
$data[206,1] = 11
$data[206,2] = 0
$data[206,3] = 10
$data[206,4] = 0
$data[206,5] = 0
$data[206,6] = 0
$data[206,7] = 1

$data[207,0] = "Seychelles"
$data[207,1] = 11
$data[207,2] = 0
$data[207,3] = 11
$data[207,4] = 0
$data[207,5] = 0
$data[207,6] = 0
$data[207,7] = 0

$data[208,0] = "Sahara Occidental"
$data[208,1] = 9
$data[208,2] = 0
$data[208,3] = 8
$data[208,4] = 0
$data[208,5] = 0
$data[208,6] = 0
$data[208,7] = 1

$data[209,0] = "Islas Virgenes Britanicas"
$data[209,1] = 8
$data[209,2] = 0
$data[209,3] = 7
$data[209,4] = 0
$data[209,5] = 0
$data[209,6] = 0
$data[209,7] = 1

$data[210,0] = "Papua Nueva Guinea"
$data[210,1] = 8
$data[210,2] = 0
$data[210,3] = 8
$data[210,4] = 0
$data[210,5] = 0
$data[210,6] = 0
$data[210,7] = 0

$data[211,0] = "Bonaire, San Eustaquio y Saba"
$data[211,1] = 7
$data[211,2] = 0
$data[211,3] = 7
$data[211,4] = 0
$data[211,5] = 0
$data[211,6] = 0
$data[211,7] = 0

$data[212,0] = "San Bartolome"
$data[212,1] = 6
$data[212,2] = 0
$data[212,3] = 6
$data[212,4] = 0
$data[212,5] = 0
$data[212,6] = 0
$data[212,7] = 0

$data[213,0] = "Lesoto"
$data[213,1] = 4
$data[213,2] = 0
$data[213,3] = 2
$data[213,4] = 2
$data[213,5] = 0
$data[213,6] = 0
$data[213,7] = 0

$data[214,0] = "Anguila"
$data[214,1] = 3
$data[214,2] = 0
$data[214,3] = 3
$data[214,4] = 0
$data[214,5] = 0
$data[214,6] = 0
$data[214,7] = 0

$data[215,0] = "San Pedro y Miquelon"
$data[215,1] = 1
$data[215,2] = 0
$data[215,3] = 1
$data[215,4] = 0
$data[215,5] = 0
$data[215,6] = 0
$data[215,7] = 0

$ws.Range("A4:H219").Value = $data

# Refresh the "last updated" timestamp banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 18 de Junio de 2020 a las 17:55"

